$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Update "Bad Drivers" table row for Intel(R) Wi-Fi 6E AX210 160MHz - 23.110.0.5
$ws.Range("C3").Value = 158
$ws.Range("D3").Value = 94.2

# Update Totals row
$ws.Range("C4").Value = 158

# Update "Good Drivers" table Total Samples for Intel(R) Wi-Fi 6E AX210 160MHz - 23.20.1.1
$ws.Range("B14").Value = 14968
